$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# New header cell D1 ("Fixed") needs the same centered style as the other header cells
$ws.Range('D1').HorizontalAlignment = -4108
$ws.Range('D1').VerticalAlignment = -4108
$ws.Range('D1').Value = 'Fixed'

# Header C1 text change
$ws.Range('C1').Value = 'Further Details'

# New row 19 B/C cells need the same centered style as the rest of column B/C
$ws.Range('B19:C19').HorizontalAlignment = -4108
$ws.Range('B19:C19').VerticalAlignment = -4108

# Row 2
$ws.Range('B2').Value = 'No IMU'
$ws.Range('C2').Value = 'missing current sensor on i2cdetect'

# Row 3
$ws.Range('B3').Value = 'broken SD card holder'
$ws.Range('C3').Value = 'needs to raspberry pi'

# Row 4
$ws.Range('B4').Value = 'water damage'
$ws.Range('C4').Value = 'only seeing two devices - missing aceel/gyro and current sensor'

# Row 5
$ws.Range('B5').Value = 'had lights off during last deployment'
$ws.Range('C5').Value = 'Seems to be running just fine now '
$ws.Range('D5').Value = 'yes'

# Row 6
$ws.Range('B6').Value = 'Water Damage and no GPS battery'
$ws.Range('C6').Value = 'seems fine now'
$ws.Range('D6').Value = 'yes'

# Row 7
$ws.Range('B7').Value = 'No IMU'
$ws.Range('C7').Value = 'Initializes IMU and powers on but never opens file to run - missing two devices so it needs some solders - missing gyro/accel and current sensor'

# Row 8
$ws.Range('B8').Value = 'not able to ssh into'
$ws.Range('C8').Value = 'no light on raspberry pi - needs new raspberry pi'

# Row 9
$ws.Range('B9').Value = 'No GPS data'
$ws.Range('C9').Value = 'not able to initialize GPS so it may need a new GPS battery '

# Row 10
$ws.Range('B10').Value = 'Lots of problems - No IMU'

# Row 11
$ws.Range('B11').Value = 'broken SD card holder'
$ws.Range('C11').Value = 'need to replace Raspberry pi'

# Row 12
$ws.Range('B12').Value = 'no none problem '
$ws.Range('C12').Value = 'ok'
$ws.Range('D12').Value = 'yes'

# Row 13
$ws.Range('B13').Value = 'needs a git pull on DUNEX branch '
$ws.Range('C13').Value = 'rebooting at minute 50 '

# Row 14
$ws.Range('B14').Value = 'I/O error on IMU '
$ws.Range('C14').Value = 'Is recroding IMU now with no problems - data all looks normal'
$ws.Range('D14').Value = 'yes'

# Row 15
$ws.Range('B15').Value = 'no none problem '
$ws.Range('C15').Value = 'looks fine'
$ws.Range('D15').Value = 'yes'

# Row 16
$ws.Range('B16').Value = 'not getting gps data'
$ws.Range('C16').Value = 'getting nan values or no data for GPS data and not software update needed since it was trying to send in a record window - old settings?'

# Row 17
$ws.Range('B17').Value = 'missing all data'
$ws.Range('C17').Value = 'needs recloning and needs service script started - no processes running '

# Row 18
$ws.Range('B18').Value = 'was missing all data'
$ws.Range('C18').Value = 'is currently getting all data'
$ws.Range('D18').Value = 'yes'

# Row 19
$ws.Range('B19').Value = 'not offloading'
$ws.Range('C19').Value = 'could have been too low of battery to actually ssh into '
$ws.Range('D19').Value = 'yes'

# Row 20
$ws.Range('B20').Value = 'water damage'

# Row 21
$ws.Range('B21').Value = 'water damage'

# Row 22
$ws.Range('B22').Value = 'water damage'

# Row 23
$ws.Range('B23').Value = 'water damage'

# Column width adjustments (ColumnWidth vs. the stored XML width differ by a
# constant padding offset in this engine, so back that out to land on the
# exact stored widths of 32 and 122.5 characters)
$ws.Columns.Item(2).ColumnWidth = 31.16666666666667
$ws.Columns.Item(3).ColumnWidth = 121.66666666666667

# Restore selection to match the author's saved cursor position
[void]$ws.Range('D27').Select()